$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7354065775871277
$ws.Range("B1").Value = 1.57258939743042
$ws.Range("C1").Value = 4.823458671569824
$ws.Range("D1").Value = 2.408813714981079
$ws.Range("E1").Value = 1.050973057746887
